$d = $word.ActiveDocument

$replacements = @(
    @("12×60=", "57×14="),
    @("98×14=", "93×57="),
    @("13×40=", "62×49="),
    @("37×84=", "89×49="),
    @("16×95=", "33×77="),
    @("35×41=", "34×98="),
    @("71×99=", "39×50="),
    @("90×46=", "18×88="),
    @("83×38=", "20×40="),
    @("44×61=", "55×58="),
    @("34×16=", "67×75="),
    @("80×83=", "93×40="),
    @("59×56=", "73×54="),
    @("66×73=", "77×50="),
    @("64×60=", "97×43="),
    @("41×80=", "52×90="),
    @("95×57=", "49×31="),
    @("90×30=", "44×35="),
    @("62×56=", "38×64="),
    @("18×27=", "64×43="),
    @("59×28=", "30×89="),
    @("13×46=", "95×61="),
    @("65×11=", "16×33="),
    @("90×85=", "67×27="),
    @("57×96=", "21×20=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
